# Auto-generated edit script applying the Rafflesia_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H17").Value = 2785
$ws.Range("J17").Value = 2785
$ws.Range("L17").Value = 8355
$ws.Range("N17").Value = -8691
$ws.Range("H70").Value = 1125
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2730
$ws.Range("H73").Value = 1125
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -2064
$ws.Range("H80").Value = 1124.75
$ws.Range("J80").Value = 749.5
$ws.Range("L80").Value = 2248.5
$ws.Range("N80").Value = -4244.5
$ws.Range("H83").Value = 1124.75
$ws.Range("J83").Value = 749.5
$ws.Range("L83").Value = 6745.5
$ws.Range("N83").Value = -16729.5
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H131").Value = 20000
$ws.Range("J131").Value = 20000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080
$ws.Range("H137").Value = 1084.5
$ws.Range("I137").Value = 926
$ws.Range("J137").Value = 1401.5
$ws.Range("K137").Value = 2778
$ws.Range("L137").Value = 4204.5
$ws.Range("M137").Value = -228
$ws.Range("N137").Value = -9304.5
$ws.Range("H138").Value = 2200
$ws.Range("I138").Value = 1000
$ws.Range("J138").Value = 2600
$ws.Range("K138").Value = 3000
$ws.Range("L138").Value = 7800
$ws.Range("M138").Value = 2140
$ws.Range("N138").Value = -18080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11064.083
$ws.Range("I32").Value = 8276.9
$ws.Range("K32").Value = 8276.9
$ws.Range("M32").Value = -7989.9
$ws.Range("H61").Value = 3192.2
$ws.Range("I61").Value = 2990.25
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2990.25
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2778.25
$ws.Range("N61").Value = -4424
$ws.Range("H97").Value = 1014.5
$ws.Range("I97").Value = 917.4
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 917.4
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -421.4
$ws.Range("N97").Value = -2492
$ws.Range("H110").Value = 848.6667
$ws.Range("I110").Value = 870.25
$ws.Range("J110").Value = 805.5
$ws.Range("K110").Value = 870.25
$ws.Range("L110").Value = 805.5
$ws.Range("M110").Value = 1174.75
$ws.Range("N110").Value = -4895.5
$ws.Range("H132").Value = 3399.25
$ws.Range("I132").Value = 2456.2856
$ws.Range("K132").Value = 7368.8568
$ws.Range("M132").Value = -4838.8568
$ws.Range("H136").Value = 3192.2
$ws.Range("I136").Value = 2990.25
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 8970.75
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -6420.75
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1248.8182
$ws.Range("I94").Value = 925.4286
$ws.Range("J94").Value = 1814.75
$ws.Range("K94").Value = 925.4286
$ws.Range("L94").Value = 1814.75
$ws.Range("M94").Value = -474.4286
$ws.Range("N94").Value = -2716.75
$ws.Range("H134").Value = 5441.909
$ws.Range("I134").Value = 1584.4
$ws.Range("J134").Value = 8656.5
$ws.Range("K134").Value = 4753.200000000001
$ws.Range("L134").Value = 25969.5
$ws.Range("M134").Value = -2218.200000000001
$ws.Range("N134").Value = -31039.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2840
$ws.Range("I2").Value = 4255
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 4255
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = -4142
$ws.Range("N2").Value = -236
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 375
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 375
$ws.Range("M22").Value = -150
$ws.Range("N22").Value = -1075
$ws.Range("H31").Value = 3758.625
$ws.Range("J31").Value = 4814.8
$ws.Range("L31").Value = 4814.8
$ws.Range("N31").Value = -5404.8
$ws.Range("H34").Value = 3758.625
$ws.Range("J34").Value = 4814.8
$ws.Range("L34").Value = 4814.8
$ws.Range("N34").Value = -5218.8
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 430.5
$ws.Range("I60").Value = 325.83334
$ws.Range("J60").Value = 587.5
$ws.Range("K60").Value = 977.5000200000001
$ws.Range("L60").Value = 1762.5
$ws.Range("M60").Value = -726.5000200000001
$ws.Range("N60").Value = -2264.5
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H131").Value = 2646.5293
$ws.Range("J131").Value = 2699.4375
$ws.Range("L131").Value = 8098.3125
$ws.Range("N131").Value = -18178.3125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 30475
$ws.Range("I80").Value = 20762.5
$ws.Range("J80").Value = 49900
$ws.Range("K80").Value = 20762.5
$ws.Range("L80").Value = 49900
$ws.Range("M80").Value = -19764.5
$ws.Range("N80").Value = -51896
$ws.Range("H83").Value = 30475
$ws.Range("I83").Value = 20762.5
$ws.Range("J83").Value = 49900
$ws.Range("K83").Value = 103812.5
$ws.Range("L83").Value = 249500
$ws.Range("M83").Value = -98820.5
$ws.Range("N83").Value = -259484
$ws.Range("H97").Value = 587.6667
$ws.Range("I97").Value = 603
$ws.Range("K97").Value = 603
$ws.Range("M97").Value = -107
$ws.Range("H126").Value = 680
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 2556
$ws.Range("I132").Value = 2074.6667
$ws.Range("K132").Value = 6224.000100000001
$ws.Range("M132").Value = -3694.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5988
$ws.Range("I132").Value = 5988
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17964
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15434
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 28000
$ws.Range("I54").Value = 28000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 28000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -27480
$ws.Range("N54").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H107").Value = 275
$ws.Range("I107").Value = 275
$ws.Range("K107").Value = 825
$ws.Range("M107").Value = 1095
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
